$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 35.69412875252057
